$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the custom message in B2 with the generic placeholder message
$ws.Range("B2").Value = "Aquí va el mensaje personalizado"

# Update the selected/active cell from B44 to B3
$ws.Range("B3").Select()
